$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# swap rows 3 and 4
$ws.Range("A3").Value2 = "2000年C"
$ws.Range("B3").Value2 = 100.7
$ws.Range("C3").Value2 = -12.6
$ws.Range("D3").Value2 = ""
$ws.Range("E3").Value2 = 787.9
$ws.Range("A4").Value2 = "2000年B"
$ws.Range("B4").Value2 = 100.6
$ws.Range("C4").Value2 = -8.6
$ws.Range("D4").Value2 = ""
$ws.Range("E4").Value2 = 511.6

# swap rows 7 and 8
$ws.Range("A7").Value2 = "2001年C"
$ws.Range("B7").Value2 = 98.7
$ws.Range("C7").Value2 = 31.7
$ws.Range("D7").Value2 = ""
$ws.Range("E7").Value2 = 865.1
$ws.Range("A8").Value2 = "2001年B"
$ws.Range("B8").Value2 = 96.7
$ws.Range("C8").Value2 = 48.1
$ws.Range("D8").Value2 = ""
$ws.Range("E8").Value2 = 567.1

# swap rows 11 and 12
$ws.Range("A11").Value2 = "2002年C"
$ws.Range("B11").Value2 = 100.6
$ws.Range("C11").Value2 = 2.4
$ws.Range("D11").Value2 = ""
$ws.Range("E11").Value2 = 966.3
$ws.Range("A12").Value2 = "2002年B"
$ws.Range("B12").Value2 = 97.8
$ws.Range("C12").Value2 = 35.9
$ws.Range("D12").Value2 = ""
$ws.Range("E12").Value2 = 624.5

# swap rows 15 and 16
$ws.Range("A15").Value2 = "2003年C"
$ws.Range("B15").Value2 = 98.90000000000001
$ws.Range("C15").Value2 = 22.6
$ws.Range("D15").Value2 = ""
$ws.Range("E15").Value2 = 1141.2
$ws.Range("A16").Value2 = "2003年B"
$ws.Range("B16").Value2 = 99.40000000000001
$ws.Range("C16").Value2 = 11.8
$ws.Range("D16").Value2 = ""
$ws.Range("E16").Value2 = 746.4

# swap rows 19 and 20
$ws.Range("A19").Value2 = "2004年C"
$ws.Range("B19").Value2 = 99.7
$ws.Range("C19").Value2 = 3.5
$ws.Range("D19").Value2 = ""
$ws.Range("E19").Value2 = 1307
$ws.Range("A20").Value2 = "2004年B"
$ws.Range("B20").Value2 = 99.5
$ws.Range("C20").Value2 = 10.9
$ws.Range("D20").Value2 = ""
$ws.Range("E20").Value2 = 878

# swap rows 23 and 24
$ws.Range("A23").Value2 = "2005年C"
$ws.Range("B23").Value2 = 98.90000000000001
$ws.Range("C23").Value2 = 14.6
$ws.Range("D23").Value2 = ""
$ws.Range("E23").Value2 = 1519.8
$ws.Range("A24").Value2 = "2005年B"
$ws.Range("B24").Value2 = 99.40000000000001
$ws.Range("C24").Value2 = 1.9
$ws.Range("D24").Value2 = ""
$ws.Range("E24").Value2 = 974.6

# swap rows 27 and 28
$ws.Range("A27").Value2 = "2006年C"
$ws.Range("B27").Value2 = 98.90000000000001
$ws.Range("C27").Value2 = 48
$ws.Range("D27").Value2 = ""
$ws.Range("E27").Value2 = 1811.6
$ws.Range("A28").Value2 = "2006年B"
$ws.Range("B28").Value2 = 99.09999999999999
$ws.Range("C28").Value2 = 24.8
$ws.Range("D28").Value2 = ""
$ws.Range("E28").Value2 = 1199

# swap rows 31 and 32
$ws.Range("A31").Value2 = "2007年C"
$ws.Range("B31").Value2 = 99.40000000000001
$ws.Range("C31").Value2 = 16.4
$ws.Range("D31").Value2 = ""
$ws.Range("E31").Value2 = 2213.2
$ws.Range("A32").Value2 = "2007年B"
$ws.Range("B32").Value2 = 99.3
$ws.Range("C32").Value2 = 15.1
$ws.Range("D32").Value2 = ""
$ws.Range("E32").Value2 = 1467.7

# swap rows 35 and 36
$ws.Range("A35").Value2 = "2008年C"
$ws.Range("B35").Value2 = 98.3
$ws.Range("C35").Value2 = 25.4
$ws.Range("D35").Value2 = ""
$ws.Range("E35").Value2 = 2332.8
$ws.Range("A36").Value2 = "2008年B"
$ws.Range("B36").Value2 = 98.2
$ws.Range("C36").Value2 = 7
$ws.Range("D36").Value2 = ""
$ws.Range("E36").Value2 = 1571.7

# swap rows 39 and 40
$ws.Range("A39").Value2 = "2009年C"
$ws.Range("B39").Value2 = 99.59999999999999
$ws.Range("C39").Value2 = 14.8
$ws.Range("D39").Value2 = ""
$ws.Range("E39").Value2 = 2508.7
$ws.Range("A40").Value2 = "2009年B"
$ws.Range("B40").Value2 = 99.7
$ws.Range("C40").Value2 = 5.8
$ws.Range("D40").Value2 = ""
$ws.Range("E40").Value2 = 1623.8

# swap rows 43 and 44
$ws.Range("A43").Value2 = "2010年C"
$ws.Range("B43").Value2 = 99.3
$ws.Range("C43").Value2 = 24.4
$ws.Range("D43").Value2 = ""
$ws.Range("E43").Value2 = 3140.5
$ws.Range("A44").Value2 = "2010年B"
$ws.Range("B44").Value2 = 99.09999999999999
$ws.Range("C44").Value2 = 29.3
$ws.Range("D44").Value2 = ""
$ws.Range("E44").Value2 = 2070.1

# swap rows 47 and 48
$ws.Range("A47").Value2 = "2011年C"
$ws.Range("B47").Value2 = 99.59999999999999
$ws.Range("C47").Value2 = 30.5
$ws.Range("D47").Value2 = ""
$ws.Range("E47").Value2 = 3429.2
$ws.Range("A48").Value2 = "2011年B"
$ws.Range("B48").Value2 = 99.59999999999999
$ws.Range("C48").Value2 = 22
$ws.Range("D48").Value2 = ""
$ws.Range("E48").Value2 = 2225.9

# swap rows 51 and 52
$ws.Range("A51").Value2 = "2012年C"
$ws.Range("B51").Value2 = 99.8
$ws.Range("C51").Value2 = 5.8
$ws.Range("D51").Value2 = ""
$ws.Range("E51").Value2 = 3818
$ws.Range("A52").Value2 = "2012年B"
$ws.Range("B52").Value2 = 99
$ws.Range("C52").Value2 = 17.9
$ws.Range("D52").Value2 = ""
$ws.Range("E52").Value2 = 2565.8

# swap rows 55 and 56
$ws.Range("A55").Value2 = "2016年C"
$ws.Range("B55").Value2 = 95.5
$ws.Range("C55").Value2 = 5.4
$ws.Range("D55").Value2 = -1.2
$ws.Range("E55").Value2 = 5714.95296
$ws.Range("A56").Value2 = "2016年B"
$ws.Range("B56").Value2 = 96.3
$ws.Range("C56").Value2 = 11.3
$ws.Range("D56").Value2 = 0.1
$ws.Range("E56").Value2 = 3798.38829

# swap rows 59 and 60
$ws.Range("A59").Value2 = "2017年C"
$ws.Range("B59").Value2 = 99.59999999999999
$ws.Range("C59").Value2 = 19.7
$ws.Range("D59").Value2 = -0.2
$ws.Range("E59").Value2 = 6348.86711
$ws.Range("A60").Value2 = "2017年B"
$ws.Range("B60").Value2 = 98.5
$ws.Range("C60").Value2 = 27.3
$ws.Range("D60").Value2 = -0.5
$ws.Range("E60").Value2 = 4160.62423

# swap rows 63 and 64
$ws.Range("A63").Value2 = "2018年C"
$ws.Range("B63").Value2 = 99.5
$ws.Range("C63").Value2 = 10.7
$ws.Range("D63").Value2 = 0.3
$ws.Range("E63").Value2 = 6238.11389
$ws.Range("A64").Value2 = "2018年B"
$ws.Range("B64").Value2 = 98.7
$ws.Range("C64").Value2 = 15.8
$ws.Range("D64").Value2 = 0.5
$ws.Range("E64").Value2 = 3993.51771

# swap rows 67 and 68
$ws.Range("A67").Value2 = "2019年C"
$ws.Range("B67").Value2 = 99.40000000000001
$ws.Range("C67").Value2 = 7.4
$ws.Range("D67").Value2 = -0.1
$ws.Range("E67").Value2 = 6942.36652
$ws.Range("A68").Value2 = "2019年B"
$ws.Range("B68").Value2 = 99.2
$ws.Range("C68").Value2 = 14.4
$ws.Range("D68").Value2 = 0.2
$ws.Range("E68").Value2 = 4457.47572

# remove columns F and G entirely
$ws.Range("F1:G69").Clear()